# Append two more daily rows (2025-12-06 and 2025-12-07) to the "Chart" sheet,
# matching the pattern of the existing rows: Date (text), Invalid count, Valid count.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Format column A for the new rows as Text first so the date-like strings are
# stored literally (matching the existing rows) instead of being auto-parsed
# into Excel date serial numbers.
$ws.Range("A63:A64").NumberFormat = "@"

$ws.Range("A63").Value = "2025-12-06"
$ws.Range("B63").Value = 0
$ws.Range("C63").Value = 25

$ws.Range("A64").Value = "2025-12-07"
$ws.Range("B64").Value = 0
$ws.Range("C64").Value = 26
